$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 97
$ws1.Range("F3").Value = 12050
$ws1.Range("F5").Value = 230
$ws1.Range("F6").Value = 365
$ws1.Range("F7").Value = 228
$ws1.Range("F8").Value = 11938
$ws1.Range("F9").Value = 502
$ws1.Range("F11").Value = 110
$ws1.Range("F12").Value = 583
$ws1.Range("F14").Value = 5910
$ws1.Range("F15").Value = 130
$ws1.Range("F16").Value = 3556
$ws1.Range("F17").Value = 198
$ws1.Range("F18").Value = 30

# ---- Sheet "演出" (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 3

# ---- Sheet "全部类型" (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 97
$ws4.Range("F4").Value = 3
$ws4.Range("F5").Value = 12050
$ws4.Range("F7").Value = 230
$ws4.Range("F9").Value = 365
$ws4.Range("F10").Value = 228
$ws4.Range("F11").Value = 11938
$ws4.Range("F12").Value = 502
$ws4.Range("F14").Value = 110
$ws4.Range("F15").Value = 583
$ws4.Range("F18").Value = 5910
$ws4.Range("F19").Value = 130
$ws4.Range("F20").Value = 3556
$ws4.Range("F21").Value = 198
$ws4.Range("F22").Value = 30
